$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 281.73
$ws.Range("D2").Value = 127632.65
$ws.Range("E2").Value = 0.62
$ws.Range("F2").Value = 188.75
$ws.Range("G2").Value = 76489.55
$ws.Range("H2").Value = 0.74
$ws.Range("I2").Value = 123.83
$ws.Range("J2").Value = 27215.84
$ws.Range("K2").Value = 0.93
